$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 82

# Column A holds the date as text (matching the existing "YYYY-MM-DD" shared
# strings), so force text interpretation to avoid Excel auto-converting the
# literal into a date serial number.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2020-08-20"
$cellA.NumberFormat = "General"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 543806
$ws.Cells.Item($row, 3).Value = 599525
$ws.Cells.Item($row, 4).Value = 82786
$ws.Cells.Item($row, 5).Value = 59106
$ws.Cells.Item($row, 6).Value = 26.02
